$wb = $excel.ActiveWorkbook

# --- 1. Update shared text labels in "BAU Emissions" (": NoSettings" -> ": test") ---
$wsBAU = $wb.Worksheets.Item("BAU Emissions")

$wsBAU.Range("A4").Value = "Process Emissions before CCS[agriculture and forestry 01T03,CO2] : test"
$wsBAU.Range("A5").Value = "Process Emissions before CCS[coal mining 05,CO2] : test"
$wsBAU.Range("A6").Value = "Process Emissions before CCS[oil and gas extraction 06,CO2] : test"
$wsBAU.Range("A7").Value = "Process Emissions before CCS[other mining and quarrying 07T08,CO2] : test"
$wsBAU.Range("A8").Value = "Process Emissions before CCS[food beverage and tobacco 10T12,CO2] : test"
$wsBAU.Range("A9").Value = "Process Emissions before CCS[textiles apparel and leather 13T15,CO2] : test"
$wsBAU.Range("A10").Value = "Process Emissions before CCS[wood products 16,CO2] : test"
$wsBAU.Range("A11").Value = "Process Emissions before CCS[pulp paper and printing 17T18,CO2] : test"
$wsBAU.Range("A12").Value = "Process Emissions before CCS[refined petroleum and coke 19,CO2] : test"
$wsBAU.Range("A13").Value = "Process Emissions before CCS[chemicals 20,CO2] : test"
$wsBAU.Range("A14").Value = "Process Emissions before CCS[rubber and plastic products 22,CO2] : test"
$wsBAU.Range("A15").Value = "Process Emissions before CCS[glass and glass products 231,CO2] : test"
$wsBAU.Range("A16").Value = "Process Emissions before CCS[cement and other nonmetallic minerals 239,CO2] : test"
$wsBAU.Range("A17").Value = "Process Emissions before CCS[iron and steel 241,CO2] : test"
$wsBAU.Range("A18").Value = "Process Emissions before CCS[other metals 242,CO2] : test"
$wsBAU.Range("A19").Value = "Process Emissions before CCS[metal products except machinery and vehicles 25,CO2] : test"
$wsBAU.Range("A20").Value = "Process Emissions before CCS[computers and electronics 26,CO2] : test"
$wsBAU.Range("A21").Value = "Process Emissions before CCS[appliances and electrical equipment 27,CO2] : test"
$wsBAU.Range("A22").Value = "Process Emissions before CCS[other machinery 28,CO2] : test"
$wsBAU.Range("A23").Value = "Process Emissions before CCS[road vehicles 29,CO2] : test"
$wsBAU.Range("A24").Value = "Process Emissions before CCS[nonroad vehicles 30,CO2] : test"
$wsBAU.Range("A25").Value = "Process Emissions before CCS[other manufacturing 31T33,CO2] : test"
$wsBAU.Range("A26").Value = "Process Emissions before CCS[energy pipelines and gas processing 352T353,CO2] : test"
$wsBAU.Range("A27").Value = "Process Emissions before CCS[water and waste 36T39,CO2] : test"
$wsBAU.Range("A28").Value = "Process Emissions before CCS[construction 41T43,CO2] : test"
$wsBAU.Range("A31").Value = "Industrial Sector Energy Related Emissions before CCS[electricity if,agriculture and forestry 01T03,CO2] : test"
$wsBAU.Range("A32").Value = "Industrial Sector Energy Related Emissions before CCS[electricity if,coal mining 05,CO2] : test"
$wsBAU.Range("A33").Value = "Industrial Sector Energy Related Emissions before CCS[electricity if,oil and gas extraction 06,CO2] : test"
$wsBAU.Range("A34").Value = "Industrial Sector Energy Related Emissions before CCS[electricity if,other mining and quarrying 07T08,CO2] : test"
$wsBAU.Range("A35").Value = "Industrial Sector Energy Related Emissions before CCS[electricity if,food beverage and tobacco 10T12,CO2] : test"
$wsBAU.Range("A36").Value = "Industrial Sector Energy Related Emissions before CCS[electricity if,textiles apparel and leather 13T15,CO2] : test"
$wsBAU.Range("A37").Value = "Industrial Sector Energy Related Emissions before CCS[electricity if,wood products 16,CO2] : test"
$wsBAU.Range("A38").Value = "Industrial Sector Energy Related Emissions before CCS[electricity if,pulp paper and printing 17T18,CO2] : test"
$wsBAU.Range("A39").Value = "Industrial Sector Energy Related Emissions before CCS[electricity if,refined petroleum and coke 19,CO2] : test"
$wsBAU.Range("A40").Value = "Industrial Sector Energy Related Emissions before CCS[electricity if,chemicals 20,CO2] : test"
$wsBAU.Range("A41").Value = "Industrial Sector Energy Related Emissions before CCS[electricity if,rubber and plastic products 22,CO2] : test"
$wsBAU.Range("A42").Value = "Industrial Sector Energy Related Emissions before CCS[electricity if,glass and glass products 231,CO2] : test"
$wsBAU.Range("A43").Value = "Industrial Sector Energy Related Emissions before CCS[electricity if,cement and other nonmetallic minerals 239,CO2] : test"
$wsBAU.Range("A44").Value = "Industrial Sector Energy Related Emissions before CCS[electricity if,iron and steel 241,CO2] : test"
$wsBAU.Range("A45").Value = "Industrial Sector Energy Related Emissions before CCS[electricity if,other metals 242,CO2] : test"
$wsBAU.Range("A46").Value = "Industrial Sector Energy Related Emissions before CCS[electricity if,metal products except machinery and vehicles 25,CO2] : test"
$wsBAU.Range("A47").Value = "Industrial Sector Energy Related Emissions before CCS[electricity if,computers and electronics 26,CO2] : test"
$wsBAU.Range("A48").Value = "Industrial Sector Energy Related Emissions before CCS[electricity if,appliances and electrical equipment 27,CO2] : test"
$wsBAU.Range("A49").Value = "Industrial Sector Energy Related Emissions before CCS[electricity if,other machinery 28,CO2] : test"
$wsBAU.Range("A50").Value = "Industrial Sector Energy Related Emissions before CCS[electricity if,road vehicles 29,CO2] : test"
$wsBAU.Range("A51").Value = "Industrial Sector Energy Related Emissions before CCS[electricity if,nonroad vehicles 30,CO2] : test"
$wsBAU.Range("A52").Value = "Industrial Sector Energy Related Emissions before CCS[electricity if,other manufacturing 31T33,CO2] : test"
$wsBAU.Range("A53").Value = "Industrial Sector Energy Related Emissions before CCS[electricity if,energy pipelines and gas processing 352T353,CO2] : test"
$wsBAU.Range("A54").Value = "Industrial Sector Energy Related Emissions before CCS[electricity if,water and waste 36T39,CO2] : test"
$wsBAU.Range("A55").Value = "Industrial Sector Energy Related Emissions before CCS[electricity if,construction 41T43,CO2] : test"
$wsBAU.Range("A56").Value = "Industrial Sector Energy Related Emissions before CCS[hard coal if,agriculture and forestry 01T03,CO2] : test"
$wsBAU.Range("A57").Value = "Industrial Sector Energy Related Emissions before CCS[hard coal if,coal mining 05,CO2] : test"
$wsBAU.Range("A58").Value = "Industrial Sector Energy Related Emissions before CCS[hard coal if,oil and gas extraction 06,CO2] : test"
$wsBAU.Range("A59").Value = "Industrial Sector Energy Related Emissions before CCS[hard coal if,other mining and quarrying 07T08,CO2] : test"
$wsBAU.Range("A60").Value = "Industrial Sector Energy Related Emissions before CCS[hard coal if,food beverage and tobacco 10T12,CO2] : test"
$wsBAU.Range("A61").Value = "Industrial Sector Energy Related Emissions before CCS[hard coal if,textiles apparel and leather 13T15,CO2] : test"
$wsBAU.Range("A62").Value = "Industrial Sector Energy Related Emissions before CCS[hard coal if,wood products 16,CO2] : test"
$wsBAU.Range("A63").Value = "Industrial Sector Energy Related Emissions before CCS[hard coal if,pulp paper and printing 17T18,CO2] : test"
$wsBAU.Range("A64").Value = "Industrial Sector Energy Related Emissions before CCS[hard coal if,refined petroleum and coke 19,CO2] : test"
$wsBAU.Range("A65").Value = "Industrial Sector Energy Related Emissions before CCS[hard coal if,chemicals 20,CO2] : test"
$wsBAU.Range("A66").Value = "Industrial Sector Energy Related Emissions before CCS[hard coal if,rubber and plastic products 22,CO2] : test"
$wsBAU.Range("A67").Value = "Industrial Sector Energy Related Emissions before CCS[hard coal if,glass and glass products 231,CO2] : test"
$wsBAU.Range("A68").Value = "Industrial Sector Energy Related Emissions before CCS[hard coal if,cement and other nonmetallic minerals 239,CO2] : test"
$wsBAU.Range("A69").Value = "Industrial Sector Energy Related Emissions before CCS[hard coal if,iron and steel 241,CO2] : test"
$wsBAU.Range("A70").Value = "Industrial Sector Energy Related Emissions before CCS[hard coal if,other metals 242,CO2] : test"
$wsBAU.Range("A71").Value = "Industrial Sector Energy Related Emissions before CCS[hard coal if,metal products except machinery and vehicles 25,CO2] : test"
$wsBAU.Range("A72").Value = "Industrial Sector Energy Related Emissions before CCS[hard coal if,computers and electronics 26,CO2] : test"
$wsBAU.Range("A73").Value = "Industrial Sector Energy Related Emissions before CCS[hard coal if,appliances and electrical equipment 27,CO2] : test"
$wsBAU.Range("A74").Value = "Industrial Sector Energy Related Emissions before CCS[hard coal if,other machinery 28,CO2] : test"
$wsBAU.Range("A75").Value = "Industrial Sector Energy Related Emissions before CCS[hard coal if,road vehicles 29,CO2] : test"
$wsBAU.Range("A76").Value = "Industrial Sector Energy Related Emissions before CCS[hard coal if,nonroad vehicles 30,CO2] : test"
$wsBAU.Range("A77").Value = "Industrial Sector Energy Related Emissions before CCS[hard coal if,other manufacturing 31T33,CO2] : test"
$wsBAU.Range("A78").Value = "Industrial Sector Energy Related Emissions before CCS[hard coal if,energy pipelines and gas processing 352T353,CO2] : test"
$wsBAU.Range("A79").Value = "Industrial Sector Energy Related Emissions before CCS[hard coal if,water and waste 36T39,CO2] : test"
$wsBAU.Range("A80").Value = "Industrial Sector Energy Related Emissions before CCS[hard coal if,construction 41T43,CO2] : test"
$wsBAU.Range("A81").Value = "Industrial Sector Energy Related Emissions before CCS[natural gas if,agriculture and forestry 01T03,CO2] : test"
$wsBAU.Range("A82").Value = "Industrial Sector Energy Related Emissions before CCS[natural gas if,coal mining 05,CO2] : test"
$wsBAU.Range("A83").Value = "Industrial Sector Energy Related Emissions before CCS[natural gas if,oil and gas extraction 06,CO2] : test"
$wsBAU.Range("A84").Value = "Industrial Sector Energy Related Emissions before CCS[natural gas if,other mining and quarrying 07T08,CO2] : test"
$wsBAU.Range("A85").Value = "Industrial Sector Energy Related Emissions before CCS[natural gas if,food beverage and tobacco 10T12,CO2] : test"
$wsBAU.Range("A86").Value = "Industrial Sector Energy Related Emissions before CCS[natural gas if,textiles apparel and leather 13T15,CO2] : test"
$wsBAU.Range("A87").Value = "Industrial Sector Energy Related Emissions before CCS[natural gas if,wood products 16,CO2] : test"
$wsBAU.Range("A88").Value = "Industrial Sector Energy Related Emissions before CCS[natural gas if,pulp paper and printing 17T18,CO2] : test"
$wsBAU.Range("A89").Value = "Industrial Sector Energy Related Emissions before CCS[natural gas if,refined petroleum and coke 19,CO2] : test"
$wsBAU.Range("A90").Value = "Industrial Sector Energy Related Emissions before CCS[natural gas if,chemicals 20,CO2] : test"
$wsBAU.Range("A91").Value = "Industrial Sector Energy Related Emissions before CCS[natural gas if,rubber and plastic products 22,CO2] : test"
$wsBAU.Range("A92").Value = "Industrial Sector Energy Related Emissions before CCS[natural gas if,glass and glass products 231,CO2] : test"
$wsBAU.Range("A93").Value = "Industrial Sector Energy Related Emissions before CCS[natural gas if,cement and other nonmetallic minerals 239,CO2] : test"
$wsBAU.Range("A94").Value = "Industrial Sector Energy Related Emissions before CCS[natural gas if,iron and steel 241,CO2] : test"
$wsBAU.Range("A95").Value = "Industrial Sector Energy Related Emissions before CCS[natural gas if,other metals 242,CO2] : test"
$wsBAU.Range("A96").Value = "Industrial Sector Energy Related Emissions before CCS[natural gas if,metal products except machinery and vehicles 25,CO2] : test"
$wsBAU.Range("A97").Value = "Industrial Sector Energy Related Emissions before CCS[natural gas if,computers and electronics 26,CO2] : test"
$wsBAU.Range("A98").Value = "Industrial Sector Energy Related Emissions before CCS[natural gas if,appliances and electrical equipment 27,CO2] : test"
$wsBAU.Range("A99").Value = "Industrial Sector Energy Related Emissions before CCS[natural gas if,other machinery 28,CO2] : test"
$wsBAU.Range("A100").Value = "Industrial Sector Energy Related Emissions before CCS[natural gas if,road vehicles 29,CO2] : test"
$wsBAU.Range("A101").Value = "Industrial Sector Energy Related Emissions before CCS[natural gas if,nonroad vehicles 30,CO2] : test"
$wsBAU.Range("A102").Value = "Industrial Sector Energy Related Emissions before CCS[natural gas if,other manufacturing 31T33,CO2] : test"
$wsBAU.Range("A103").Value = "Industrial Sector Energy Related Emissions before CCS[natural gas if,energy pipelines and gas processing 352T353,CO2] : test"
$wsBAU.Range("A104").Value = "Industrial Sector Energy Related Emissions before CCS[natural gas if,water and waste 36T39,CO2] : test"
$wsBAU.Range("A105").Value = "Industrial Sector Energy Related Emissions before CCS[natural gas if,construction 41T43,CO2] : test"
$wsBAU.Range("A106").Value = "Industrial Sector Energy Related Emissions before CCS[biomass if,agriculture and forestry 01T03,CO2] : test"
$wsBAU.Range("A107").Value = "Industrial Sector Energy Related Emissions before CCS[biomass if,coal mining 05,CO2] : test"
$wsBAU.Range("A108").Value = "Industrial Sector Energy Related Emissions before CCS[biomass if,oil and gas extraction 06,CO2] : test"
$wsBAU.Range("A109").Value = "Industrial Sector Energy Related Emissions before CCS[biomass if,other mining and quarrying 07T08,CO2] : test"
$wsBAU.Range("A110").Value = "Industrial Sector Energy Related Emissions before CCS[biomass if,food beverage and tobacco 10T12,CO2] : test"
$wsBAU.Range("A111").Value = "Industrial Sector Energy Related Emissions before CCS[biomass if,textiles apparel and leather 13T15,CO2] : test"
$wsBAU.Range("A112").Value = "Industrial Sector Energy Related Emissions before CCS[biomass if,wood products 16,CO2] : test"
$wsBAU.Range("A113").Value = "Industrial Sector Energy Related Emissions before CCS[biomass if,pulp paper and printing 17T18,CO2] : test"
$wsBAU.Range("A114").Value = "Industrial Sector Energy Related Emissions before CCS[biomass if,refined petroleum and coke 19,CO2] : test"
$wsBAU.Range("A115").Value = "Industrial Sector Energy Related Emissions before CCS[biomass if,chemicals 20,CO2] : test"
$wsBAU.Range("A116").Value = "Industrial Sector Energy Related Emissions before CCS[biomass if,rubber and plastic products 22,CO2] : test"
$wsBAU.Range("A117").Value = "Industrial Sector Energy Related Emissions before CCS[biomass if,glass and glass products 231,CO2] : test"
$wsBAU.Range("A118").Value = "Industrial Sector Energy Related Emissions before CCS[biomass if,cement and other nonmetallic minerals 239,CO2] : test"
$wsBAU.Range("A119").Value = "Industrial Sector Energy Related Emissions before CCS[biomass if,iron and steel 241,CO2] : test"
$wsBAU.Range("A120").Value = "Industrial Sector Energy Related Emissions before CCS[biomass if,other metals 242,CO2] : test"
$wsBAU.Range("A121").Value = "Industrial Sector Energy Related Emissions before CCS[biomass if,metal products except machinery and vehicles 25,CO2] : test"
$wsBAU.Range("A122").Value = "Industrial Sector Energy Related Emissions before CCS[biomass if,computers and electronics 26,CO2] : test"
$wsBAU.Range("A123").Value = "Industrial Sector Energy Related Emissions before CCS[biomass if,appliances and electrical equipment 27,CO2] : test"
$wsBAU.Range("A124").Value = "Industrial Sector Energy Related Emissions before CCS[biomass if,other machinery 28,CO2] : test"
$wsBAU.Range("A125").Value = "Industrial Sector Energy Related Emissions before CCS[biomass if,road vehicles 29,CO2] : test"
$wsBAU.Range("A126").Value = "Industrial Sector Energy Related Emissions before CCS[biomass if,nonroad vehicles 30,CO2] : test"
$wsBAU.Range("A127").Value = "Industrial Sector Energy Related Emissions before CCS[biomass if,other manufacturing 31T33,CO2] : test"
$wsBAU.Range("A128").Value = "Industrial Sector Energy Related Emissions before CCS[biomass if,energy pipelines and gas processing 352T353,CO2] : test"
$wsBAU.Range("A129").Value = "Industrial Sector Energy Related Emissions before CCS[biomass if,water and waste 36T39,CO2] : test"
$wsBAU.Range("A130").Value = "Industrial Sector Energy Related Emissions before CCS[biomass if,construction 41T43,CO2] : test"
$wsBAU.Range("A131").Value = "Industrial Sector Energy Related Emissions before CCS[petroleum diesel if,agriculture and forestry 01T03,CO2] : test"
$wsBAU.Range("A132").Value = "Industrial Sector Energy Related Emissions before CCS[petroleum diesel if,coal mining 05,CO2] : test"
$wsBAU.Range("A133").Value = "Industrial Sector Energy Related Emissions before CCS[petroleum diesel if,oil and gas extraction 06,CO2] : test"
$wsBAU.Range("A134").Value = "Industrial Sector Energy Related Emissions before CCS[petroleum diesel if,other mining and quarrying 07T08,CO2] : test"
$wsBAU.Range("A135").Value = "Industrial Sector Energy Related Emissions before CCS[petroleum diesel if,food beverage and tobacco 10T12,CO2] : test"
$wsBAU.Range("A136").Value = "Industrial Sector Energy Related Emissions before CCS[petroleum diesel if,textiles apparel and leather 13T15,CO2] : test"
$wsBAU.Range("A137").Value = "Industrial Sector Energy Related Emissions before CCS[petroleum diesel if,wood products 16,CO2] : test"
$wsBAU.Range("A138").Value = "Industrial Sector Energy Related Emissions before CCS[petroleum diesel if,pulp paper and printing 17T18,CO2] : test"
$wsBAU.Range("A139").Value = "Industrial Sector Energy Related Emissions before CCS[petroleum diesel if,refined petroleum and coke 19,CO2] : test"
$wsBAU.Range("A140").Value = "Industrial Sector Energy Related Emissions before CCS[petroleum diesel if,chemicals 20,CO2] : test"
$wsBAU.Range("A141").Value = "Industrial Sector Energy Related Emissions before CCS[petroleum diesel if,rubber and plastic products 22,CO2] : test"
$wsBAU.Range("A142").Value = "Industrial Sector Energy Related Emissions before CCS[petroleum diesel if,glass and glass products 231,CO2] : test"
$wsBAU.Range("A143").Value = "Industrial Sector Energy Related Emissions before CCS[petroleum diesel if,cement and other nonmetallic minerals 239,CO2] : test"
$wsBAU.Range("A144").Value = "Industrial Sector Energy Related Emissions before CCS[petroleum diesel if,iron and steel 241,CO2] : test"
$wsBAU.Range("A145").Value = "Industrial Sector Energy Related Emissions before CCS[petroleum diesel if,other metals 242,CO2] : test"
$wsBAU.Range("A146").Value = "Industrial Sector Energy Related Emissions before CCS[petroleum diesel if,metal products except machinery and vehicles 25,CO2] : test"
$wsBAU.Range("A147").Value = "Industrial Sector Energy Related Emissions before CCS[petroleum diesel if,computers and electronics 26,CO2] : test"
$wsBAU.Range("A148").Value = "Industrial Sector Energy Related Emissions before CCS[petroleum diesel if,appliances and electrical equipment 27,CO2] : test"
$wsBAU.Range("A149").Value = "Industrial Sector Energy Related Emissions before CCS[petroleum diesel if,other machinery 28,CO2] : test"
$wsBAU.Range("A150").Value = "Industrial Sector Energy Related Emissions before CCS[petroleum diesel if,road vehicles 29,CO2] : test"
$wsBAU.Range("A151").Value = "Industrial Sector Energy Related Emissions before CCS[petroleum diesel if,nonroad vehicles 30,CO2] : test"
$wsBAU.Range("A152").Value = "Industrial Sector Energy Related Emissions before CCS[petroleum diesel if,other manufacturing 31T33,CO2] : test"
$wsBAU.Range("A153").Value = "Industrial Sector Energy Related Emissions before CCS[petroleum diesel if,energy pipelines and gas processing 352T353,CO2] : test"
$wsBAU.Range("A154").Value = "Industrial Sector Energy Related Emissions before CCS[petroleum diesel if,water and waste 36T39,CO2] : test"
$wsBAU.Range("A155").Value = "Industrial Sector Energy Related Emissions before CCS[petroleum diesel if,construction 41T43,CO2] : test"
$wsBAU.Range("A156").Value = "Industrial Sector Energy Related Emissions before CCS[heat if,agriculture and forestry 01T03,CO2] : test"
$wsBAU.Range("A157").Value = "Industrial Sector Energy Related Emissions before CCS[heat if,coal mining 05,CO2] : test"
$wsBAU.Range("A158").Value = "Industrial Sector Energy Related Emissions before CCS[heat if,oil and gas extraction 06,CO2] : test"
$wsBAU.Range("A159").Value = "Industrial Sector Energy Related Emissions before CCS[heat if,other mining and quarrying 07T08,CO2] : test"
$wsBAU.Range("A160").Value = "Industrial Sector Energy Related Emissions before CCS[heat if,food beverage and tobacco 10T12,CO2] : test"
$wsBAU.Range("A161").Value = "Industrial Sector Energy Related Emissions before CCS[heat if,textiles apparel and leather 13T15,CO2] : test"
$wsBAU.Range("A162").Value = "Industrial Sector Energy Related Emissions before CCS[heat if,wood products 16,CO2] : test"
$wsBAU.Range("A163").Value = "Industrial Sector Energy Related Emissions before CCS[heat if,pulp paper and printing 17T18,CO2] : test"
$wsBAU.Range("A164").Value = "Industrial Sector Energy Related Emissions before CCS[heat if,refined petroleum and coke 19,CO2] : test"
$wsBAU.Range("A165").Value = "Industrial Sector Energy Related Emissions before CCS[heat if,chemicals 20,CO2] : test"
$wsBAU.Range("A166").Value = "Industrial Sector Energy Related Emissions before CCS[heat if,rubber and plastic products 22,CO2] : test"
$wsBAU.Range("A167").Value = "Industrial Sector Energy Related Emissions before CCS[heat if,glass and glass products 231,CO2] : test"
$wsBAU.Range("A168").Value = "Industrial Sector Energy Related Emissions before CCS[heat if,cement and other nonmetallic minerals 239,CO2] : test"
$wsBAU.Range("A169").Value = "Industrial Sector Energy Related Emissions before CCS[heat if,iron and steel 241,CO2] : test"
$wsBAU.Range("A170").Value = "Industrial Sector Energy Related Emissions before CCS[heat if,other metals 242,CO2] : test"
$wsBAU.Range("A171").Value = "Industrial Sector Energy Related Emissions before CCS[heat if,metal products except machinery and vehicles 25,CO2] : test"
$wsBAU.Range("A172").Value = "Industrial Sector Energy Related Emissions before CCS[heat if,computers and electronics 26,CO2] : test"
$wsBAU.Range("A173").Value = "Industrial Sector Energy Related Emissions before CCS[heat if,appliances and electrical equipment 27,CO2] : test"
$wsBAU.Range("A174").Value = "Industrial Sector Energy Related Emissions before CCS[heat if,other machinery 28,CO2] : test"
$wsBAU.Range("A175").Value = "Industrial Sector Energy Related Emissions before CCS[heat if,road vehicles 29,CO2] : test"
$wsBAU.Range("A176").Value = "Industrial Sector Energy Related Emissions before CCS[heat if,nonroad vehicles 30,CO2] : test"
$wsBAU.Range("A177").Value = "Industrial Sector Energy Related Emissions before CCS[heat if,other manufacturing 31T33,CO2] : test"
$wsBAU.Range("A178").Value = "Industrial Sector Energy Related Emissions before CCS[heat if,energy pipelines and gas processing 352T353,CO2] : test"
$wsBAU.Range("A179").Value = "Industrial Sector Energy Related Emissions before CCS[heat if,water and waste 36T39,CO2] : test"
$wsBAU.Range("A180").Value = "Industrial Sector Energy Related Emissions before CCS[heat if,construction 41T43,CO2] : test"
$wsBAU.Range("A181").Value = "Industrial Sector Energy Related Emissions before CCS[crude oil if,agriculture and forestry 01T03,CO2] : test"
$wsBAU.Range("A182").Value = "Industrial Sector Energy Related Emissions before CCS[crude oil if,coal mining 05,CO2] : test"
$wsBAU.Range("A183").Value = "Industrial Sector Energy Related Emissions before CCS[crude oil if,oil and gas extraction 06,CO2] : test"
$wsBAU.Range("A184").Value = "Industrial Sector Energy Related Emissions before CCS[crude oil if,other mining and quarrying 07T08,CO2] : test"
$wsBAU.Range("A185").Value = "Industrial Sector Energy Related Emissions before CCS[crude oil if,food beverage and tobacco 10T12,CO2] : test"
$wsBAU.Range("A186").Value = "Industrial Sector Energy Related Emissions before CCS[crude oil if,textiles apparel and leather 13T15,CO2] : test"
$wsBAU.Range("A187").Value = "Industrial Sector Energy Related Emissions before CCS[crude oil if,wood products 16,CO2] : test"
$wsBAU.Range("A188").Value = "Industrial Sector Energy Related Emissions before CCS[crude oil if,pulp paper and printing 17T18,CO2] : test"
$wsBAU.Range("A189").Value = "Industrial Sector Energy Related Emissions before CCS[crude oil if,refined petroleum and coke 19,CO2] : test"
$wsBAU.Range("A190").Value = "Industrial Sector Energy Related Emissions before CCS[crude oil if,chemicals 20,CO2] : test"
$wsBAU.Range("A191").Value = "Industrial Sector Energy Related Emissions before CCS[crude oil if,rubber and plastic products 22,CO2] : test"
$wsBAU.Range("A192").Value = "Industrial Sector Energy Related Emissions before CCS[crude oil if,glass and glass products 231,CO2] : test"
$wsBAU.Range("A193").Value = "Industrial Sector Energy Related Emissions before CCS[crude oil if,cement and other nonmetallic minerals 239,CO2] : test"
$wsBAU.Range("A194").Value = "Industrial Sector Energy Related Emissions before CCS[crude oil if,iron and steel 241,CO2] : test"
$wsBAU.Range("A195").Value = "Industrial Sector Energy Related Emissions before CCS[crude oil if,other metals 242,CO2] : test"
$wsBAU.Range("A196").Value = "Industrial Sector Energy Related Emissions before CCS[crude oil if,metal products except machinery and vehicles 25,CO2] : test"
$wsBAU.Range("A197").Value = "Industrial Sector Energy Related Emissions before CCS[crude oil if,computers and electronics 26,CO2] : test"
$wsBAU.Range("A198").Value = "Industrial Sector Energy Related Emissions before CCS[crude oil if,appliances and electrical equipment 27,CO2] : test"
$wsBAU.Range("A199").Value = "Industrial Sector Energy Related Emissions before CCS[crude oil if,other machinery 28,CO2] : test"
$wsBAU.Range("A200").Value = "Industrial Sector Energy Related Emissions before CCS[crude oil if,road vehicles 29,CO2] : test"
$wsBAU.Range("A201").Value = "Industrial Sector Energy Related Emissions before CCS[crude oil if,nonroad vehicles 30,CO2] : test"
$wsBAU.Range("A202").Value = "Industrial Sector Energy Related Emissions before CCS[crude oil if,other manufacturing 31T33,CO2] : test"
$wsBAU.Range("A203").Value = "Industrial Sector Energy Related Emissions before CCS[crude oil if,energy pipelines and gas processing 352T353,CO2] : test"
$wsBAU.Range("A204").Value = "Industrial Sector Energy Related Emissions before CCS[crude oil if,water and waste 36T39,CO2] : test"
$wsBAU.Range("A205").Value = "Industrial Sector Energy Related Emissions before CCS[crude oil if,construction 41T43,CO2] : test"
$wsBAU.Range("A206").Value = "Industrial Sector Energy Related Emissions before CCS[heavy or residual fuel oil if,agriculture and forestry 01T03,CO2] : test"
$wsBAU.Range("A207").Value = "Industrial Sector Energy Related Emissions before CCS[heavy or residual fuel oil if,coal mining 05,CO2] : test"
$wsBAU.Range("A208").Value = "Industrial Sector Energy Related Emissions before CCS[heavy or residual fuel oil if,oil and gas extraction 06,CO2] : test"
$wsBAU.Range("A209").Value = "Industrial Sector Energy Related Emissions before CCS[heavy or residual fuel oil if,other mining and quarrying 07T08,CO2] : test"
$wsBAU.Range("A210").Value = "Industrial Sector Energy Related Emissions before CCS[heavy or residual fuel oil if,food beverage and tobacco 10T12,CO2] : test"
$wsBAU.Range("A211").Value = "Industrial Sector Energy Related Emissions before CCS[heavy or residual fuel oil if,textiles apparel and leather 13T15,CO2] : test"
$wsBAU.Range("A212").Value = "Industrial Sector Energy Related Emissions before CCS[heavy or residual fuel oil if,wood products 16,CO2] : test"
$wsBAU.Range("A213").Value = "Industrial Sector Energy Related Emissions before CCS[heavy or residual fuel oil if,pulp paper and printing 17T18,CO2] : test"
$wsBAU.Range("A214").Value = "Industrial Sector Energy Related Emissions before CCS[heavy or residual fuel oil if,refined petroleum and coke 19,CO2] : test"
$wsBAU.Range("A215").Value = "Industrial Sector Energy Related Emissions before CCS[heavy or residual fuel oil if,chemicals 20,CO2] : test"
$wsBAU.Range("A216").Value = "Industrial Sector Energy Related Emissions before CCS[heavy or residual fuel oil if,rubber and plastic products 22,CO2] : test"
$wsBAU.Range("A217").Value = "Industrial Sector Energy Related Emissions before CCS[heavy or residual fuel oil if,glass and glass products 231,CO2] : test"
$wsBAU.Range("A218").Value = "Industrial Sector Energy Related Emissions before CCS[heavy or residual fuel oil if,cement and other nonmetallic minerals 239,CO2] : test"
$wsBAU.Range("A219").Value = "Industrial Sector Energy Related Emissions before CCS[heavy or residual fuel oil if,iron and steel 241,CO2] : test"
$wsBAU.Range("A220").Value = "Industrial Sector Energy Related Emissions before CCS[heavy or residual fuel oil if,other metals 242,CO2] : test"
$wsBAU.Range("A221").Value = "Industrial Sector Energy Related Emissions before CCS[heavy or residual fuel oil if,metal products except machinery and vehicles 25,CO2] : test"
$wsBAU.Range("A222").Value = "Industrial Sector Energy Related Emissions before CCS[heavy or residual fuel oil if,computers and electronics 26,CO2] : test"
$wsBAU.Range("A223").Value = "Industrial Sector Energy Related Emissions before CCS[heavy or residual fuel oil if,appliances and electrical equipment 27,CO2] : test"
$wsBAU.Range("A224").Value = "Industrial Sector Energy Related Emissions before CCS[heavy or residual fuel oil if,other machinery 28,CO2] : test"
$wsBAU.Range("A225").Value = "Industrial Sector Energy Related Emissions before CCS[heavy or residual fuel oil if,road vehicles 29,CO2] : test"
$wsBAU.Range("A226").Value = "Industrial Sector Energy Related Emissions before CCS[heavy or residual fuel oil if,nonroad vehicles 30,CO2] : test"
$wsBAU.Range("A227").Value = "Industrial Sector Energy Related Emissions before CCS[heavy or residual fuel oil if,other manufacturing 31T33,CO2] : test"
$wsBAU.Range("A228").Value = "Industrial Sector Energy Related Emissions before CCS[heavy or residual fuel oil if,energy pipelines and gas processing 352T353,CO2] : test"
$wsBAU.Range("A229").Value = "Industrial Sector Energy Related Emissions before CCS[heavy or residual fuel oil if,water and waste 36T39,CO2] : test"
$wsBAU.Range("A230").Value = "Industrial Sector Energy Related Emissions before CCS[heavy or residual fuel oil if,construction 41T43,CO2] : test"
$wsBAU.Range("A231").Value = "Industrial Sector Energy Related Emissions before CCS[LPG propane or butane if,agriculture and forestry 01T03,CO2] : test"
$wsBAU.Range("A232").Value = "Industrial Sector Energy Related Emissions before CCS[LPG propane or butane if,coal mining 05,CO2] : test"
$wsBAU.Range("A233").Value = "Industrial Sector Energy Related Emissions before CCS[LPG propane or butane if,oil and gas extraction 06,CO2] : test"
$wsBAU.Range("A234").Value = "Industrial Sector Energy Related Emissions before CCS[LPG propane or butane if,other mining and quarrying 07T08,CO2] : test"
$wsBAU.Range("A235").Value = "Industrial Sector Energy Related Emissions before CCS[LPG propane or butane if,food beverage and tobacco 10T12,CO2] : test"
$wsBAU.Range("A236").Value = "Industrial Sector Energy Related Emissions before CCS[LPG propane or butane if,textiles apparel and leather 13T15,CO2] : test"
$wsBAU.Range("A237").Value = "Industrial Sector Energy Related Emissions before CCS[LPG propane or butane if,wood products 16,CO2] : test"
$wsBAU.Range("A238").Value = "Industrial Sector Energy Related Emissions before CCS[LPG propane or butane if,pulp paper and printing 17T18,CO2] : test"
$wsBAU.Range("A239").Value = "Industrial Sector Energy Related Emissions before CCS[LPG propane or butane if,refined petroleum and coke 19,CO2] : test"
$wsBAU.Range("A240").Value = "Industrial Sector Energy Related Emissions before CCS[LPG propane or butane if,chemicals 20,CO2] : test"
$wsBAU.Range("A241").Value = "Industrial Sector Energy Related Emissions before CCS[LPG propane or butane if,rubber and plastic products 22,CO2] : test"
$wsBAU.Range("A242").Value = "Industrial Sector Energy Related Emissions before CCS[LPG propane or butane if,glass and glass products 231,CO2] : test"
$wsBAU.Range("A243").Value = "Industrial Sector Energy Related Emissions before CCS[LPG propane or butane if,cement and other nonmetallic minerals 239,CO2] : test"
$wsBAU.Range("A244").Value = "Industrial Sector Energy Related Emissions before CCS[LPG propane or butane if,iron and steel 241,CO2] : test"
$wsBAU.Range("A245").Value = "Industrial Sector Energy Related Emissions before CCS[LPG propane or butane if,other metals 242,CO2] : test"
$wsBAU.Range("A246").Value = "Industrial Sector Energy Related Emissions before CCS[LPG propane or butane if,metal products except machinery and vehicles 25,CO2] : test"
$wsBAU.Range("A247").Value = "Industrial Sector Energy Related Emissions before CCS[LPG propane or butane if,computers and electronics 26,CO2] : test"
$wsBAU.Range("A248").Value = "Industrial Sector Energy Related Emissions before CCS[LPG propane or butane if,appliances and electrical equipment 27,CO2] : test"
$wsBAU.Range("A249").Value = "Industrial Sector Energy Related Emissions before CCS[LPG propane or butane if,other machinery 28,CO2] : test"
$wsBAU.Range("A250").Value = "Industrial Sector Energy Related Emissions before CCS[LPG propane or butane if,road vehicles 29,CO2] : test"
$wsBAU.Range("A251").Value = "Industrial Sector Energy Related Emissions before CCS[LPG propane or butane if,nonroad vehicles 30,CO2] : test"
$wsBAU.Range("A252").Value = "Industrial Sector Energy Related Emissions before CCS[LPG propane or butane if,other manufacturing 31T33,CO2] : test"
$wsBAU.Range("A253").Value = "Industrial Sector Energy Related Emissions before CCS[LPG propane or butane if,energy pipelines and gas processing 352T353,CO2] : test"
$wsBAU.Range("A254").Value = "Industrial Sector Energy Related Emissions before CCS[LPG propane or butane if,water and waste 36T39,CO2] : test"
$wsBAU.Range("A255").Value = "Industrial Sector Energy Related Emissions before CCS[LPG propane or butane if,construction 41T43,CO2] : test"
$wsBAU.Range("A256").Value = "Industrial Sector Energy Related Emissions before CCS[hydrogen if,agriculture and forestry 01T03,CO2] : test"
$wsBAU.Range("A257").Value = "Industrial Sector Energy Related Emissions before CCS[hydrogen if,coal mining 05,CO2] : test"
$wsBAU.Range("A258").Value = "Industrial Sector Energy Related Emissions before CCS[hydrogen if,oil and gas extraction 06,CO2] : test"
$wsBAU.Range("A259").Value = "Industrial Sector Energy Related Emissions before CCS[hydrogen if,other mining and quarrying 07T08,CO2] : test"
$wsBAU.Range("A260").Value = "Industrial Sector Energy Related Emissions before CCS[hydrogen if,food beverage and tobacco 10T12,CO2] : test"
$wsBAU.Range("A261").Value = "Industrial Sector Energy Related Emissions before CCS[hydrogen if,textiles apparel and leather 13T15,CO2] : test"
$wsBAU.Range("A262").Value = "Industrial Sector Energy Related Emissions before CCS[hydrogen if,wood products 16,CO2] : test"
$wsBAU.Range("A263").Value = "Industrial Sector Energy Related Emissions before CCS[hydrogen if,pulp paper and printing 17T18,CO2] : test"
$wsBAU.Range("A264").Value = "Industrial Sector Energy Related Emissions before CCS[hydrogen if,refined petroleum and coke 19,CO2] : test"
$wsBAU.Range("A265").Value = "Industrial Sector Energy Related Emissions before CCS[hydrogen if,chemicals 20,CO2] : test"
$wsBAU.Range("A266").Value = "Industrial Sector Energy Related Emissions before CCS[hydrogen if,rubber and plastic products 22,CO2] : test"
$wsBAU.Range("A267").Value = "Industrial Sector Energy Related Emissions before CCS[hydrogen if,glass and glass products 231,CO2] : test"
$wsBAU.Range("A268").Value = "Industrial Sector Energy Related Emissions before CCS[hydrogen if,cement and other nonmetallic minerals 239,CO2] : test"
$wsBAU.Range("A269").Value = "Industrial Sector Energy Related Emissions before CCS[hydrogen if,iron and steel 241,CO2] : test"
$wsBAU.Range("A270").Value = "Industrial Sector Energy Related Emissions before CCS[hydrogen if,other metals 242,CO2] : test"
$wsBAU.Range("A271").Value = "Industrial Sector Energy Related Emissions before CCS[hydrogen if,metal products except machinery and vehicles 25,CO2] : test"
$wsBAU.Range("A272").Value = "Industrial Sector Energy Related Emissions before CCS[hydrogen if,computers and electronics 26,CO2] : test"
$wsBAU.Range("A273").Value = "Industrial Sector Energy Related Emissions before CCS[hydrogen if,appliances and electrical equipment 27,CO2] : test"
$wsBAU.Range("A274").Value = "Industrial Sector Energy Related Emissions before CCS[hydrogen if,other machinery 28,CO2] : test"
$wsBAU.Range("A275").Value = "Industrial Sector Energy Related Emissions before CCS[hydrogen if,road vehicles 29,CO2] : test"
$wsBAU.Range("A276").Value = "Industrial Sector Energy Related Emissions before CCS[hydrogen if,nonroad vehicles 30,CO2] : test"
$wsBAU.Range("A277").Value = "Industrial Sector Energy Related Emissions before CCS[hydrogen if,other manufacturing 31T33,CO2] : test"
$wsBAU.Range("A278").Value = "Industrial Sector Energy Related Emissions before CCS[hydrogen if,energy pipelines and gas processing 352T353,CO2] : test"
$wsBAU.Range("A279").Value = "Industrial Sector Energy Related Emissions before CCS[hydrogen if,water and waste 36T39,CO2] : test"
$wsBAU.Range("A280").Value = "Industrial Sector Energy Related Emissions before CCS[hydrogen if,construction 41T43,CO2] : test"

# --- 2. Update the "last updated" date on the About sheet (C1: 45369 -> 45387) ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45387

# --- 3. Update the refreshed projection values in row 94 (columns M:AE) ---
$wsBAU.Range("M94").Value = 1001080
$wsBAU.Range("N94").Value = 2002150
$wsBAU.Range("O94").Value = 3003230
$wsBAU.Range("P94").Value = 4004300
$wsBAU.Range("Q94").Value = 5005380
$wsBAU.Range("R94").Value = 5005380
$wsBAU.Range("S94").Value = 5005380
$wsBAU.Range("T94").Value = 5005380
$wsBAU.Range("U94").Value = 5005380
$wsBAU.Range("V94").Value = 5005380
$wsBAU.Range("W94").Value = 5005380
$wsBAU.Range("X94").Value = 5005380
$wsBAU.Range("Y94").Value = 5005380
$wsBAU.Range("Z94").Value = 5005380
$wsBAU.Range("AA94").Value = 5005380
$wsBAU.Range("AB94").Value = 5005380
$wsBAU.Range("AC94").Value = 5005380
$wsBAU.Range("AD94").Value = 5005380
$wsBAU.Range("AE94").Value = 5005380

# --- 4. Update sheet view / active-tab state to match the saved workbook ---
# "BAU Emissions": selection becomes the block A30:AE280 (active cell A30)
$wsBAU.Activate()
$wsBAU.Range("A30:AE280").Select()

# "Current and Planned Capacity": keep its prior selection, just no longer the active tab
$wsCurrent = $wb.Worksheets.Item("Current and Planned Capacity")
$wsCurrent.Activate()
$wsCurrent.Range("I59").Select()

# "About" ends up the active tab (tabSelected moves here; workbook opens on this sheet)
$wsAbout.Activate()
$wsAbout.Range("E29").Select()

